# Update countries & provincias Spain
# Applies the COVID-19 data refresh captured in the diff:
#   - bump the "Datos actualizados" timestamp (15:35 -> 16:05)
#   - refresh the numeric stats for a handful of countries (rows keep their
#     position because they are sorted by total cases, column B)
#   - Mali's case count overtook Hong Kong's, so Mali is inserted above
#     Hong Kong/Haiti in the (descending, by total cases) list; Hong Kong and
#     Haiti's existing rows simply slide down one row with their data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Timestamp banner in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 26 de Mayo de 2020 a las 16:05"

# 2. Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 1708597
$ws.Cells.Item(4, 3).Value = 2371
$ws.Cells.Item(4, 4).Value = 465668
$ws.Cells.Item(4, 5).Value = 1143082
$ws.Cells.Item(4, 7).Value = 42
$ws.Cells.Item(4, 8).Value = 99847

# 3. Reino Unido (row 8)
$ws.Cells.Item(8, 7).Value = 134
$ws.Cells.Item(8, 8).Value = 37048

# 4. India (row 13)
$ws.Cells.Item(13, 2).Value = 147144
$ws.Cells.Item(13, 3).Value = 2194
$ws.Cells.Item(13, 4).Value = 61923
$ws.Cells.Item(13, 5).Value = 81024
$ws.Cells.Item(13, 7).Value = 25
$ws.Cells.Item(13, 8).Value = 4197

# 5. Mali moves ahead of Hong Kong / Haiti (rows 107-109), pushing the
#    latter two down a row with their previous totals untouched.
$ws.Cells.Item(107, 1).Value = "Mali"
$ws.Cells.Item(107, 2).Value = 1077
$ws.Cells.Item(107, 3).Value = 18
$ws.Cells.Item(107, 4).Value = 617
$ws.Cells.Item(107, 5).Value = 390
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 3
$ws.Cells.Item(107, 8).Value = 70

$ws.Cells.Item(108, 1).Value = "Hong Kong"
$ws.Cells.Item(108, 2).Value = 1066
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = 1033
$ws.Cells.Item(108, 5).Value = 29
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 4

$ws.Cells.Item(109, 1).Value = "Haiti"
$ws.Cells.Item(109, 2).Value = 1063
$ws.Cells.Item(109, 3).Value = 105
$ws.Cells.Item(109, 4).Value = 22
$ws.Cells.Item(109, 5).Value = 1010
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 4
$ws.Cells.Item(109, 8).Value = 31

# 6. Birmania (row 157)
$ws.Cells.Item(157, 2).Value = 206
$ws.Cells.Item(157, 3).Value = 3
$ws.Cells.Item(157, 4).Value = 124
$ws.Cells.Item(157, 5).Value = 76
